$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of row 51 column A (bold/center/border) onto new A-column cells
$srcStyleCell = $ws.Range("A51")

# Row 52
$ws.Range("A52").Value = 50
$srcStyleCell.Copy() | Out-Null
$ws.Range("A52").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B52").Value = 'Die Bell 525 hat einen ersten Flug mit SAF absolviert                            
premium'
$ws.Range("D52").Value = 'Nachdem Bell Textron seit März 2021 nachhaltig produzierten Treibstoff (SAF) für seine eigene Vorführ- und Trainingsflotte nutzt, hat der Hersteller nun bekanntgegeben, […]'
$ws.Range("E52").Value = 'https://aerobuzz.de/helikopter/die-bell-525-hat-einen-ersten-flug-mit-saf-absolviert/'
$ws.Range("F52").Value = $true
$ws.Range("G52").Value = '2021.11.25 - 15:58'

# Row 53
$ws.Range("A53").Value = 51
$srcStyleCell.Copy() | Out-Null
$ws.Range("A53").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B53").Value = '["['' - ETC Awarded Contract from the University of North Dakota’s John D. Odegard School of Aerospace Sciences'', "]'
$ws.Range("E53").Value = 'https://www.etcusa.com/etc-awarded-contract-from-the-university-of-north-dakotas-john-d-odegard-school-of-aerospace-sciences/'
$ws.Range("F53").Value = $false
$ws.Range("G53").Value = '2021.11.25 - 15:58'

# Row 54
$ws.Range("A54").Value = 52
$srcStyleCell.Copy() | Out-Null
$ws.Range("A54").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B54").Value = '["['' - ETC Announces Notice of Annual Meeting of Shareholders'', "]'
$ws.Range("E54").Value = 'https://www.etcusa.com/etc-announces-notice-of-annual-meeting-of-shareholders-2/'
$ws.Range("F54").Value = $false
$ws.Range("G54").Value = '2021.11.25 - 15:58'

# Row 55
$ws.Range("A55").Value = 53
$srcStyleCell.Copy() | Out-Null
$ws.Range("A55").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B55").Value = '["['' ETC Announces Fiscal 2022 Second Quarter Results'', "]'
$ws.Range("E55").Value = 'https://www.etcusa.com/etc-announces-fiscal-2022-second-quarter-results/'
$ws.Range("F55").Value = $false
$ws.Range("G55").Value = '2021.11.25 - 15:58'

# Row 56
$ws.Range("A56").Value = 54
$srcStyleCell.Copy() | Out-Null
$ws.Range("A56").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B56").Value = 'Reiser Simulation and Training teams up with RS Flight Systems to provide new product line on FNPTs and FTDs to the helicpoter training market'
$ws.Range("C56").Value = 'Berg, Germany, 30.09.2021 Reiser Simulation and Training GmbH (RST) has teamed up with RS Flight Systems GmbH (RSFS) to provide innovative and cost-effective FNPT and FTD products to the market. The new product line called “F-light line” for “Flightsimulator light” will feature the well-proven Level D capable flight model together with the Helionix® avionics replica, […]'
$ws.Range("E56").Value = 'https://www.reiser-st.com/reiser-simulation-and-training-teams-up-with-rs-flight-systems-to-provide-new-product-line-on-fnpts-andftds-to-the-helicpoter-training-market/'
$ws.Range("F56").Value = $false
$ws.Range("G56").Value = '2021.11.25 - 15:58'

# Row 57
$ws.Range("A57").Value = 55
$srcStyleCell.Copy() | Out-Null
$ws.Range("A57").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B57").Value = 'Reiser Simulation and Training GmbH to Expand Helicopter Full Flight Simulator Portfolio'
$ws.Range("C57").Value = 'Berg, Germany, September 16, 2021 Reiser Simulation and Training GmbH (RST) expands its portfolio of innovative flight training solutions. The German simulator manufacturer is proud of having welcomed Swiss Air-Rescue (Rega) and Leonardo Helicopters to its headquarters for the signing of their latest contract. RST has been awarded a contract from Rega for the provision […]'
$ws.Range("E57").Value = 'https://www.reiser-st.com/reiser-simulation-and-training-gmbh-to-expand-helicopter-full-flight-simulator-portfolio/'
$ws.Range("F57").Value = $false
$ws.Range("G57").Value = '2021.11.25 - 15:58'

# Row 58
$ws.Range("A58").Value = 56
$srcStyleCell.Copy() | Out-Null
$ws.Range("A58").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B58").Value = 'Successful EASA Level D Qualification<br>of Reiser Simulation and Training GmbH (RST) H135 Full Flight Simulator<br>at Lufthansa Aviation Training GmbH (LAT)'
$ws.Range("C58").Value = 'Berg / Frankfurt a. Main / Germany, May 31, 2021 RST’s newest Airbus H135 full flight simulator was initially qualified to the highest standard Level D by the German Federal Aviation Office (Luftfahrt-Bundesamt / LBA) according to European Aviation Safety Agency (EASA) regulations. The most modern FFS helicopter type H135 is thus ready for training […]'
$ws.Range("E58").Value = 'https://www.reiser-st.com/successful-easa-level-d-qualificationof-reiser-simulation-and-training-gmbh-rst-h135-full-flight-simulatorat-lufthansa-aviation-training-gmbh-lat/'
$ws.Range("F58").Value = $false
$ws.Range("G58").Value = '2021.11.25 - 15:58'

# Row 59
$ws.Range("A59").Value = 57
$srcStyleCell.Copy() | Out-Null
$ws.Range("A59").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B59").Value = '
RHT/UWETS Asia Repeat Order '
$ws.Range("C59").Value = '
				News				
2021-02-15 
'
$ws.Range("D59").Value = '
			We are proud to announce a repeat order for a combined Rescue Hoist and Under Water Escape Training System. The...		'
$ws.Range("E59").Value = 'https://www.amst.co.at/news-aerospace-medicine/rht-uwets-asia-repeat-order/'
$ws.Range("F59").Value = $false
$ws.Range("G59").Value = '2021.11.25 - 15:58'

# Row 60
$ws.Range("A60").Value = 58
$srcStyleCell.Copy() | Out-Null
$ws.Range("A60").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B60").Value = '
Normobaric Hypoxia Acceptance '
$ws.Range("C60").Value = '
				News				
2021-02-09 
'
$ws.Range("D60").Value = '
			AMST successfully handed over two Reduced Oxygen Breathing Devices (ROBDs) for normobaric hypoxia training to two of its AIRFOX customers....		'
$ws.Range("E60").Value = 'https://www.amst.co.at/news-aerospace-medicine/normobaric-hypoxia-acceptance/'
$ws.Range("F60").Value = $false
$ws.Range("G60").Value = '2021.11.25 - 15:58'

# Row 61
$ws.Range("A61").Value = 59
$srcStyleCell.Copy() | Out-Null
$ws.Range("A61").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B61").Value = '
Charity Donation '
$ws.Range("C61").Value = '
				News				
2020-12-16 
'
$ws.Range("D61").Value = '
			Already in early autumn, the staff at our headquarters in Ranshofen decided to donate to a charitable cause instead of...		'
$ws.Range("E61").Value = 'https://www.amst.co.at/news-aerospace-medicine/charity-donation/'
$ws.Range("F61").Value = $false
$ws.Range("G61").Value = '2021.11.25 - 15:58'

# Row 62
$ws.Range("A62").Value = 60
$srcStyleCell.Copy() | Out-Null
$ws.Range("A62").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B62").Value = '
AMST at WATS 2021 '
$ws.Range("C62").Value = '
				News				
2021-06-11 
'
$ws.Range("D62").Value = '
			AMST is participating in the World Aviation Training Summit (WATS) taking place on 15 and 16 of June 2021 in Orlando Florida. Join us at our booth #223 to have...		'
$ws.Range("E62").Value = 'https://www.amst.co.at/news-civil-aviation/amst-at-wats-2021/'
$ws.Range("F62").Value = $false
$ws.Range("G62").Value = '2021.11.25 - 15:58'

# Row 63
$ws.Range("A63").Value = 61
$srcStyleCell.Copy() | Out-Null
$ws.Range("A63").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B63").Value = '
FFS UPRT Upgrade '
$ws.Range("C63").Value = '
				News				
2021-02-15 
'
$ws.Range("D63").Value = '
			Upset Prevention and Recovery Training (UPRT) remains a crucial topic in our industry. AMST has developed and delivered a range of solutions to this demanding requirement. We have applied almost...		'
$ws.Range("E63").Value = 'https://www.amst.co.at/news-civil-aviation/ffs-uprt-upgrade/'
$ws.Range("F63").Value = $false
$ws.Range("G63").Value = '2021.11.25 - 15:58'

# Row 64
$ws.Range("A64").Value = 62
$srcStyleCell.Copy() | Out-Null
$ws.Range("A64").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B64").Value = '
Engl Flightteam Cooperation '
$ws.Range("C64").Value = '
				News				
2020-11-20 
'
$ws.Range("D64").Value = '
			Engl Flightteam, a flight training provider, signed a partnership agreement with AMST in October 2020. The cooperation focusses on Spatial Disorientation training on our AIRFOX ASD and Multi Crew Coordination...		'
$ws.Range("E64").Value = 'https://www.amst.co.at/news-aerospace-medicine/engl-flightteam-cooperation/'
$ws.Range("F64").Value = $false
$ws.Range("G64").Value = '2021.11.25 - 15:58'

Write-Host "Rows 52-64 added. UsedRange:" $ws.UsedRange.Address()
